# "Added data for ProvarCache"
# The Provar automation tool appends a fresh block of RMA numbers/ids each
# time it runs against the "RMA Details Maintenance Grid" sheet, and then
# repoints the three data rows (rows 2-4) at the newest block (here the
# "RMA-N3P5" block). Reproduce that: refresh the live RMA numbers / line ids
# for rows 2-4 with the new cached values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 (RMA-N3P5-001 / line 1-1)
$ws.Range("E2").Value = "RMA-N3P5-001"
$ws.Range("F2").Value = "RMA-N3P5-1-1"
$ws.Range("J2").Value = "a7s5f000000xKM0AAM"

# Row 3 (RMA-N3P5-002 / line 1-2)
$ws.Range("E3").Value = "RMA-N3P5-002"
$ws.Range("F3").Value = "RMA-N3P5-1-2"
$ws.Range("J3").Value = "a7s5f000000xKM1AAM"

# Row 4 (RMA-N3P5-003 / line 1-3)
$ws.Range("E4").Value = "RMA-N3P5-003"
$ws.Range("F4").Value = "RMA-N3P5-1-3"
$ws.Range("J4").Value = "a7s5f000000xKM2AAM"

# The longer RMA ids widen the best-fit columns F (RMA line) and J (SF id)
# slightly, matching the sheet's existing bestFit/customWidth columns.
$ws.Columns(6).ColumnWidth = 13.666667
$ws.Columns(10).ColumnWidth = 21.166667
